# HCAP-1233: update the seeded test fixture value and cell selection
# on Sheet1 of participant_post_hire_status.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data fix: B2 changes from 68 to 26
$ws.Range("B2").Value = 26

# Reflect the new active selection (was B3, now B2)
$ws.Range("B2").Select()
